# Auto-generated edit script applying the odds-data refresh diff
# (row swaps for corrected match ordering + odds updates for upcoming fixtures)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23
$ws.Range("B23").Value = 5965129
$ws.Range("C23").Value = "Costa Rica Primera Division"
$ws.Range("D23").Value = "Costa Rica Primera Division"
$ws.Range("E23").Value = 44955.875
$ws.Range("F23").Value = "Alajuelense"
$ws.Range("G23").Value = "Sporting San Jose"
$ws.Range("H23").Value = 3
$ws.Range("I23").Value = 1
$ws.Range("J23").Value = "H"
$ws.Range("K23").Value = 1.727
$ws.Range("L23").Value = 3.5
$ws.Range("M23").Value = 3.8
$ws.Range("N23").Value = 1.615
$ws.Range("O23").Value = 3.6
$ws.Range("P23").Value = 4.333
$ws.Range("Q23").Value = -0.75
$ws.Range("R23").Value = 1.875
$ws.Range("S23").Value = 1.925
$ws.Range("T23").Value = 2.75
$ws.Range("U23").Value = 1.85
$ws.Range("V23").Value = 1.95
$ws.Range("W23").Value = 0.615
$ws.Range("X23").Value = -1
$ws.Range("Y23").Value = -1
$ws.Range("Z23").Value = 0.875
$ws.Range("AA23").Value = -1
$ws.Range("AB23").Value = 0.8500000000000001
$ws.Range("AC23").Value = -1

# Row 24
$ws.Range("B24").Value = 5965131
$ws.Range("C24").Value = "Costa Rica Primera Division"
$ws.Range("D24").Value = "Costa Rica Primera Division"
$ws.Range("E24").Value = 44955.875
$ws.Range("F24").Value = "Guadalupe FC"
$ws.Range("G24").Value = "Cartagines"
$ws.Range("H24").Value = 5
$ws.Range("I24").Value = 1
$ws.Range("J24").Value = "H"
$ws.Range("K24").Value = 2.375
$ws.Range("L24").Value = 3.25
$ws.Range("M24").Value = 2.6
$ws.Range("N24").Value = 2.7
$ws.Range("O24").Value = 3.2
$ws.Range("P24").Value = 2.3
$ws.Range("Q24").Value = 0
$ws.Range("R24").Value = 2.05
$ws.Range("S24").Value = 1.75
$ws.Range("T24").Value = 2.5
$ws.Range("U24").Value = 1.85
$ws.Range("V24").Value = 1.95
$ws.Range("W24").Value = 1.7
$ws.Range("X24").Value = -1
$ws.Range("Y24").Value = -1
$ws.Range("Z24").Value = 1.05
$ws.Range("AA24").Value = -1
$ws.Range("AB24").Value = 0.8500000000000001
$ws.Range("AC24").Value = -1

# Row 90
$ws.Range("B90").Value = 5965176
$ws.Range("C90").Value = "Costa Rica Primera Division"
$ws.Range("D90").Value = "Costa Rica Primera Division"
$ws.Range("E90").Value = 45015.95833333334
$ws.Range("F90").Value = "AD San Carlos"
$ws.Range("G90").Value = "AD Grecia"
$ws.Range("H90").Value = 3
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = "H"
$ws.Range("K90").Value = 1.571
$ws.Range("L90").Value = 4
$ws.Range("M90").Value = 4.75
$ws.Range("N90").Value = 1.363
$ws.Range("O90").Value = 5
$ws.Range("P90").Value = 6
$ws.Range("Q90").Value = -1.25
$ws.Range("R90").Value = 1.85
$ws.Range("S90").Value = 1.95
$ws.Range("T90").Value = 3
$ws.Range("U90").Value = 1.925
$ws.Range("V90").Value = 1.875
$ws.Range("W90").Value = 0.363
$ws.Range("X90").Value = -1
$ws.Range("Y90").Value = -1
$ws.Range("Z90").Value = 0.8500000000000001
$ws.Range("AA90").Value = -1
$ws.Range("AB90").Value = 0
$ws.Range("AC90").Value = -0

# Row 91
$ws.Range("B91").Value = 5965174
$ws.Range("C91").Value = "Costa Rica Primera Division"
$ws.Range("D91").Value = "Costa Rica Primera Division"
$ws.Range("E91").Value = 45015.95833333334
$ws.Range("F91").Value = "Cartagines"
$ws.Range("G91").Value = "Guadalupe FC"
$ws.Range("H91").Value = 3
$ws.Range("I91").Value = 2
$ws.Range("J91").Value = "H"
$ws.Range("K91").Value = 1.571
$ws.Range("L91").Value = 4
$ws.Range("M91").Value = 5
$ws.Range("N91").Value = 1.363
$ws.Range("O91").Value = 5
$ws.Range("P91").Value = 6.5
$ws.Range("Q91").Value = -1.5
$ws.Range("R91").Value = 2.025
$ws.Range("S91").Value = 1.775
$ws.Range("T91").Value = 3.25
$ws.Range("U91").Value = 1.925
$ws.Range("V91").Value = 1.875
$ws.Range("W91").Value = 0.363
$ws.Range("X91").Value = -1
$ws.Range("Y91").Value = -1
$ws.Range("Z91").Value = -1
$ws.Range("AA91").Value = 0.7749999999999999
$ws.Range("AB91").Value = 0.925
$ws.Range("AC91").Value = -1

# Row 132
$ws.Range("B132").Value = 5965203
$ws.Range("C132").Value = "Costa Rica Primera Division"
$ws.Range("D132").Value = "Costa Rica Primera Division"
$ws.Range("E132").Value = 45053.75
$ws.Range("F132").Value = "Cartagines"
$ws.Range("G132").Value = "Sporting San Jose"
$ws.Range("H132").Value = 3
$ws.Range("I132").Value = 2
$ws.Range("J132").Value = "H"
$ws.Range("K132").Value = 1.85
$ws.Range("L132").Value = 3.5
$ws.Range("M132").Value = 3.4
$ws.Range("N132").Value = 2.2
$ws.Range("O132").Value = 3.25
$ws.Range("P132").Value = 2.875
$ws.Range("Q132").Value = -0.25
$ws.Range("R132").Value = 1.975
$ws.Range("S132").Value = 1.825
$ws.Range("T132").Value = 2.5
$ws.Range("U132").Value = 1.925
$ws.Range("V132").Value = 1.875
$ws.Range("W132").Value = 1.2
$ws.Range("X132").Value = -1
$ws.Range("Y132").Value = -1
$ws.Range("Z132").Value = 0.9750000000000001
$ws.Range("AA132").Value = -1
$ws.Range("AB132").Value = 0.925
$ws.Range("AC132").Value = -1

# Row 133
$ws.Range("B133").Value = 5965205
$ws.Range("C133").Value = "Costa Rica Primera Division"
$ws.Range("D133").Value = "Costa Rica Primera Division"
$ws.Range("E133").Value = 45053.75
$ws.Range("F133").Value = "Puntarenas"
$ws.Range("G133").Value = "Herediano"
$ws.Range("H133").Value = 1
$ws.Range("I133").Value = 2
$ws.Range("J133").Value = "A"
$ws.Range("K133").Value = 3.5
$ws.Range("L133").Value = 3.3
$ws.Range("M133").Value = 1.909
$ws.Range("N133").Value = 4.5
$ws.Range("O133").Value = 3.6
$ws.Range("P133").Value = 1.65
$ws.Range("Q133").Value = 0.75
$ws.Range("R133").Value = 1.925
$ws.Range("S133").Value = 1.875
$ws.Range("T133").Value = 2.5
$ws.Range("U133").Value = 1.9
$ws.Range("V133").Value = 1.9
$ws.Range("W133").Value = -1
$ws.Range("X133").Value = -1
$ws.Range("Y133").Value = 0.6499999999999999
$ws.Range("Z133").Value = -0.5
$ws.Range("AA133").Value = 0.4375
$ws.Range("AB133").Value = 0.8999999999999999
$ws.Range("AC133").Value = -1

# Row 231
$ws.Range("B231").Value = 6782568
$ws.Range("C231").Value = "Costa Rica Primera Division"
$ws.Range("D231").Value = "Costa Rica Primera Division"
$ws.Range("E231").Value = 45220.83333333334
$ws.Range("F231").Value = "Sporting San Jose"
$ws.Range("G231").Value = "AD Guanacasteca"
$ws.Range("H231").Value = 1
$ws.Range("I231").Value = 1
$ws.Range("J231").Value = "D"
$ws.Range("K231").Value = 1.909
$ws.Range("L231").Value = 3.6
$ws.Range("M231").Value = 3.3
$ws.Range("N231").Value = 2
$ws.Range("O231").Value = 3.6
$ws.Range("P231").Value = 3.1
$ws.Range("Q231").Value = -0.5
$ws.Range("R231").Value = 2
$ws.Range("S231").Value = 1.8
$ws.Range("T231").Value = 2.5
$ws.Range("U231").Value = 1.825
$ws.Range("V231").Value = 1.975
$ws.Range("W231").Value = -1
$ws.Range("X231").Value = 2.6
$ws.Range("Y231").Value = -1
$ws.Range("Z231").Value = -1
$ws.Range("AA231").Value = 0.8
$ws.Range("AB231").Value = -1
$ws.Range("AC231").Value = 0.9750000000000001

# Row 232
$ws.Range("B232").Value = 6782566
$ws.Range("C232").Value = "Costa Rica Primera Division"
$ws.Range("D232").Value = "Costa Rica Primera Division"
$ws.Range("E232").Value = 45220.83333333334
$ws.Range("F232").Value = "Cartagines"
$ws.Range("G232").Value = "Deportivo Saprissa"
$ws.Range("H232").Value = 0
$ws.Range("I232").Value = 4
$ws.Range("J232").Value = "A"
$ws.Range("K232").Value = 3.2
$ws.Range("L232").Value = 3.4
$ws.Range("M232").Value = 2
$ws.Range("N232").Value = 2.9
$ws.Range("O232").Value = 3.5
$ws.Range("P232").Value = 2.15
$ws.Range("Q232").Value = 0.25
$ws.Range("R232").Value = 1.875
$ws.Range("S232").Value = 1.925
$ws.Range("T232").Value = 3
$ws.Range("U232").Value = 1.975
$ws.Range("V232").Value = 1.825
$ws.Range("W232").Value = -1
$ws.Range("X232").Value = -1
$ws.Range("Y232").Value = 1.15
$ws.Range("Z232").Value = -1
$ws.Range("AA232").Value = 0.925
$ws.Range("AB232").Value = 0.9750000000000001
$ws.Range("AC232").Value = -1

# Row 250
$ws.Range("B250").Value = 6782581
$ws.Range("C250").Value = "Costa Rica Primera Division"
$ws.Range("D250").Value = "Costa Rica Primera Division"
$ws.Range("E250").Value = 45238.875
$ws.Range("F250").Value = "Alajuelense"
$ws.Range("G250").Value = "AD Grecia"
$ws.Range("H250").Value = 2
$ws.Range("I250").Value = 0
$ws.Range("J250").Value = "H"
$ws.Range("K250").Value = 1.181
$ws.Range("L250").Value = 6.5
$ws.Range("M250").Value = 11
$ws.Range("N250").Value = 1.25
$ws.Range("O250").Value = 5
$ws.Range("P250").Value = 9
$ws.Range("Q250").Value = -1.75
$ws.Range("R250").Value = 1.975
$ws.Range("S250").Value = 1.825
$ws.Range("T250").Value = 3.25
$ws.Range("U250").Value = 2
$ws.Range("V250").Value = 1.8
$ws.Range("W250").Value = 0.25
$ws.Range("X250").Value = -1
$ws.Range("Y250").Value = -1
$ws.Range("Z250").Value = 0.4875
$ws.Range("AA250").Value = -0.5
$ws.Range("AB250").Value = -1
$ws.Range("AC250").Value = 0.8

# Row 251
$ws.Range("B251").Value = 6782579
$ws.Range("C251").Value = "Costa Rica Primera Division"
$ws.Range("D251").Value = "Costa Rica Primera Division"
$ws.Range("E251").Value = 45238.875
$ws.Range("F251").Value = "Santos de Gupiles"
$ws.Range("G251").Value = "AD San Carlos"
$ws.Range("H251").Value = 0
$ws.Range("I251").Value = 2
$ws.Range("J251").Value = "A"
$ws.Range("K251").Value = 2.4
$ws.Range("L251").Value = 3.3
$ws.Range("M251").Value = 2.7
$ws.Range("N251").Value = 2.375
$ws.Range("O251").Value = 3.4
$ws.Range("P251").Value = 2.8
$ws.Range("Q251").Value = -0.25
$ws.Range("R251").Value = 2
$ws.Range("S251").Value = 1.8
$ws.Range("T251").Value = 2.5
$ws.Range("U251").Value = 1.875
$ws.Range("V251").Value = 1.925
$ws.Range("W251").Value = -1
$ws.Range("X251").Value = -1
$ws.Range("Y251").Value = 1.8
$ws.Range("Z251").Value = -1
$ws.Range("AA251").Value = 0.8
$ws.Range("AB251").Value = -1
$ws.Range("AC251").Value = 0.925

# Row 269
$ws.Range("B269").Value = 6782598
$ws.Range("C269").Value = "Costa Rica Primera Division"
$ws.Range("D269").Value = "Costa Rica Primera Division"
$ws.Range("E269").Value = 45255.95833333334
$ws.Range("F269").Value = "Municipal Perez Zeledon"
$ws.Range("G269").Value = "Cartagines"
$ws.Range("H269").Value = 1
$ws.Range("I269").Value = 0
$ws.Range("J269").Value = "H"
$ws.Range("K269").Value = 4.5
$ws.Range("L269").Value = 3.75
$ws.Range("M269").Value = 1.615
$ws.Range("N269").Value = 3.4
$ws.Range("O269").Value = 3.4
$ws.Range("P269").Value = 1.85
$ws.Range("Q269").Value = 0.5
$ws.Range("R269").Value = 1.8
$ws.Range("S269").Value = 2
$ws.Range("T269").Value = 2.75
$ws.Range("U269").Value = 1.9
$ws.Range("V269").Value = 1.9
$ws.Range("W269").Value = 2.4
$ws.Range("X269").Value = -1
$ws.Range("Y269").Value = -1
$ws.Range("Z269").Value = 0.8
$ws.Range("AA269").Value = -1
$ws.Range("AB269").Value = -1
$ws.Range("AC269").Value = 0.8999999999999999

# Row 271
$ws.Range("B271").Value = 6782595
$ws.Range("C271").Value = "Costa Rica Primera Division"
$ws.Range("D271").Value = "Costa Rica Primera Division"
$ws.Range("E271").Value = 45255.95833333334
$ws.Range("F271").Value = "Herediano"
$ws.Range("G271").Value = "Sporting San Jose"
$ws.Range("H271").Value = 3
$ws.Range("I271").Value = 0
$ws.Range("J271").Value = "H"
$ws.Range("K271").Value = 1.4
$ws.Range("L271").Value = 4.75
$ws.Range("M271").Value = 7
$ws.Range("N271").Value = 1.363
$ws.Range("O271").Value = 4.75
$ws.Range("P271").Value = 8.5
$ws.Range("Q271").Value = -1.25
$ws.Range("R271").Value = 1.8
$ws.Range("S271").Value = 2
$ws.Range("T271").Value = 3
$ws.Range("U271").Value = 1.95
$ws.Range("V271").Value = 1.85
$ws.Range("W271").Value = 0.363
$ws.Range("X271").Value = -1
$ws.Range("Y271").Value = -1
$ws.Range("Z271").Value = 0.8
$ws.Range("AA271").Value = -1
$ws.Range("AB271").Value = 0
$ws.Range("AC271").Value = -0

# Row 329
$ws.Range("N329").Value = 4.333
$ws.Range("O329").Value = 3.3
$ws.Range("P329").Value = 1.75
$ws.Range("Q329").Value = 0.5
$ws.Range("R329").Value = 2.025
$ws.Range("S329").Value = 1.775
$ws.Range("U329").Value = 1.8
$ws.Range("V329").Value = 2

# Row 330
$ws.Range("N330").Value = 1.833
$ws.Range("O330").Value = 3.6
$ws.Range("P330").Value = 3.75
$ws.Range("R330").Value = 1.85
$ws.Range("S330").Value = 1.95
$ws.Range("U330").Value = 1.875
$ws.Range("V330").Value = 1.925

# Row 331
$ws.Range("E331").Value = 45346.9375
$ws.Range("N331").Value = 1.615
$ws.Range("O331").Value = 3.6
$ws.Range("P331").Value = 4.75
$ws.Range("Q331").Value = -0.75
$ws.Range("R331").Value = 1.775
$ws.Range("S331").Value = 2.025
$ws.Range("U331").Value = 1.825
$ws.Range("V331").Value = 1.975

# Row 332
$ws.Range("B332").Value = 7623916
$ws.Range("C332").Value = "Costa Rica Primera Division"
$ws.Range("D332").Value = "Costa Rica Primera Division"
$ws.Range("E332").Value = 45347.75
$ws.Range("F332").Value = "Santos de Gupiles"
$ws.Range("G332").Value = "AD Grecia"
$ws.Range("K332").Value = 2.05
$ws.Range("L332").Value = 3.3
$ws.Range("M332").Value = 3.2
$ws.Range("N332").Value = 2
$ws.Range("O332").Value = 3.3
$ws.Range("P332").Value = 3.4
$ws.Range("Q332").Value = -0.5
$ws.Range("R332").Value = 2.025
$ws.Range("S332").Value = 1.775
$ws.Range("T332").Value = 2.5
$ws.Range("U332").Value = 1.825
$ws.Range("V332").Value = 1.975
$ws.Range("W332").Value = 0
$ws.Range("X332").Value = 0
$ws.Range("Y332").Value = 0
$ws.Range("Z332").Value = 0
$ws.Range("AA332").Value = 0

# Row 333
$ws.Range("B333").Value = 7623919
$ws.Range("C333").Value = "Costa Rica Primera Division"
$ws.Range("D333").Value = "Costa Rica Primera Division"
$ws.Range("E333").Value = 45347.75
$ws.Range("F333").Value = "Municipal Liberia"
$ws.Range("G333").Value = "Sporting San Jose"
$ws.Range("K333").Value = 1.75
$ws.Range("L333").Value = 3.6
$ws.Range("M333").Value = 3.8
$ws.Range("N333").Value = 1.75
$ws.Range("O333").Value = 3.6
$ws.Range("P333").Value = 3.75
$ws.Range("Q333").Value = -0.75
$ws.Range("R333").Value = 2.025
$ws.Range("S333").Value = 1.775
$ws.Range("T333").Value = 2.75
$ws.Range("U333").Value = 1.975
$ws.Range("V333").Value = 1.825
$ws.Range("W333").Value = 0
$ws.Range("X333").Value = 0
$ws.Range("Y333").Value = 0
$ws.Range("Z333").Value = 0
$ws.Range("AA333").Value = 0
